$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E9").Value = 17.39270000000002
$ws.Range("E18").Value = 17.88870000000002
$ws.Range("E20").Value = 16.1097
$ws.Range("E27").Value = 16.71669999999999
$ws.Range("E35").Value = 16.1121
$ws.Range("E69").Value = 17.36270000000002
$ws.Range("E76").Value = 16.12819999999999
$ws.Range("E78").Value = 16.59070000000003
$ws.Range("E82").Value = 16.658
$ws.Range("E83").Value = 16.5601
$ws.Range("E93").Value = 18.12090000000002
